$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.611.94'
$ws.Range("E2").Value = '  +0.94%  '

# Row 3
$ws.Range("D3").Value = '3.386.50'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.40'
$ws.Range("E5").Value = '  +0.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.95'
$ws.Range("E6").Value = '  +0.78%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = '3.384.76'
$ws.Range("E8").Value = '  +0.17%  '

# Row 9
$ws.Range("E9").Value = '  -0.62%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.48'
$ws.Range("E10").Value = '  -1.21%  '

# Row 11
$ws.Range("E11").Value = '  +2.09%  '

# Row 12
$ws.Range("E12").Value = '  +0.77%  '

# Row 13
$ws.Range("D13").Value = '3.959.91'
$ws.Range("E13").Value = '  -0.06%  '

# Row 14
$ws.Range("E14").Value = '  +1.77%  '

# Row 15
$ws.Range("E15").Value = '  +1.56%  '

# Row 16
$ws.Range("D16").Value = '3.386.01'
$ws.Range("E16").Value = '  -0.19%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.71'
$ws.Range("E17").Value = '  +2.29%  '

# Row 18
$ws.Range("D18").Value = '61.702.90'
$ws.Range("E18").Value = '  +0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.20'
$ws.Range("E19").Value = '  +1.38%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.51'
$ws.Range("E20").Value = '  +1.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("E21").Value = '  +0.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '378.78'
$ws.Range("E22").Value = '  +1.38%  '

# Row 23
$ws.Range("E23").Value = '  -1.13%  '

# Row 24
$ws.Range("D24").Value = '3.523.46'
$ws.Range("E24").Value = '  +0.12%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("E26").Value = '  +6.35%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.21'
$ws.Range("E27").Value = '  +0.90%  '

# Row 28
$ws.Range("E28").Value = '  +4.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").Value = '  -1.21%  '

# Row 30
$ws.Range("E30").Value = '  +0.05%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.19'
$ws.Range("E31").Value = '  +1.02%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.160'
$ws.Range("E32").Value = '  +2.98%  '

# Row 33
$ws.Range("E33").Value = '  +0.94%  '

# Row 34
$ws.Range("E34").Value = '  +0.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.38'
$ws.Range("E35").Value = '  -0.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.36'
$ws.Range("E36").Value = '  -3.58%  '

# Row 37
$ws.Range("E37").Value = '  -0.38%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.86'
$ws.Range("E38").Value = '  -0.81%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.78'
$ws.Range("E39").Value = '  +1.10%  '

# Row 40
$ws.Range("E40").Value = '  -0.38%  '

# Row 41
$ws.Range("E41").Value = '  +2.80%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.13%  '

# Row 43
$ws.Range("E43").Value = '  +2.36%  '

# Row 44
$ws.Range("E44").Value = '  +7.58%  '

# Row 45
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.91'
$ws.Range("E46").Value = '  +6.70%  '

# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.33'
$ws.Range("E47").Value = '  -0.06%  '

# Row 48
$ws.Range("E48").Value = '  -1.72%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.79'
$ws.Range("E49").Value = '  -0.76%  '

# Row 50
$ws.Range("D50").Value = '2.336.38'
$ws.Range("E50").Value = '  +5.82%  '

# Row 51
$ws.Range("E51").Value = '  +1.44%  '
